# Updated cryptos list on Tue Jul  2 16:35:58 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for each
# coin row on the active worksheet, and fixes the ordering of the last two rows
# (dogwifhat / FirstDigitalUSD were swapped along with their Price/Volume data).
#
# All Price/Volume cells in this sheet are stored as text (not numbers), so any
# new value that Excel would otherwise auto-convert to a number (e.g. "0.999",
# "148.16") is written with a leading apostrophe to force a text cell, exactly
# like typing it into Excel by hand would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '61.879.27'; ForceText = $false },
    @{ Cell = 'E2'; Value = '  -2.19%  '; ForceText = $false },
    @{ Cell = 'D3'; Value = '3.411.61'; ForceText = $false },
    @{ Cell = 'E3'; Value = '  -1.71%  '; ForceText = $false },
    @{ Cell = 'D4'; Value = '0.999'; ForceText = $true },
    @{ Cell = 'E4'; Value = '  +0.05%  '; ForceText = $false },
    @{ Cell = 'D5'; Value = '574.88'; ForceText = $true },
    @{ Cell = 'E5'; Value = '  -0.80%  '; ForceText = $false },
    @{ Cell = 'D6'; Value = '148.16'; ForceText = $true },
    @{ Cell = 'E6'; Value = '  -0.14%  '; ForceText = $false },
    @{ Cell = 'E7'; Value = '  +0.12%  '; ForceText = $false },
    @{ Cell = 'E8'; Value = '  +0.54%  '; ForceText = $false },
    @{ Cell = 'E9'; Value = '  +4.78%  '; ForceText = $false },
    @{ Cell = 'E10'; Value = '  -2.04%  '; ForceText = $false },
    @{ Cell = 'E11'; Value = '  +1.81%  '; ForceText = $false },
    @{ Cell = 'D12'; Value = '3.994.48'; ForceText = $false },
    @{ Cell = 'E12'; Value = '  -1.66%  '; ForceText = $false },
    @{ Cell = 'E13'; Value = '  +0.03%  '; ForceText = $false },
    @{ Cell = 'D14'; Value = '28.18'; ForceText = $true },
    @{ Cell = 'E14'; Value = '  -5.46%  '; ForceText = $false },
    @{ Cell = 'D15'; Value = '3.415.73'; ForceText = $false },
    @{ Cell = 'E15'; Value = '  -1.82%  '; ForceText = $false },
    @{ Cell = 'E16'; Value = '  -1.17%  '; ForceText = $false },
    @{ Cell = 'D17'; Value = '61.880.59'; ForceText = $false },
    @{ Cell = 'E17'; Value = '  -2.03%  '; ForceText = $false },
    @{ Cell = 'E18'; Value = '  +0.27%  '; ForceText = $false },
    @{ Cell = 'D19'; Value = '14.47'; ForceText = $true },
    @{ Cell = 'E19'; Value = '  +0.18%  '; ForceText = $false },
    @{ Cell = 'D20'; Value = '8.96'; ForceText = $true },
    @{ Cell = 'E20'; Value = '  -3.43%  '; ForceText = $false },
    @{ Cell = 'D21'; Value = '380.76'; ForceText = $true },
    @{ Cell = 'E21'; Value = '  -2.19%  '; ForceText = $false },
    @{ Cell = 'E22'; Value = '  +1.23%  '; ForceText = $false },
    @{ Cell = 'D23'; Value = '75.05'; ForceText = $true },
    @{ Cell = 'E23'; Value = '  +0.42%  '; ForceText = $false },
    @{ Cell = 'E24'; Value = '  +0.02%  '; ForceText = $false },
    @{ Cell = 'D25'; Value = '3.559.12'; ForceText = $false },
    @{ Cell = 'E25'; Value = '  -1.43%  '; ForceText = $false },
    @{ Cell = 'E26'; Value = '  -4.57%  '; ForceText = $false },
    @{ Cell = 'E27'; Value = '  +0.36%  '; ForceText = $false },
    @{ Cell = 'D28'; Value = '7.61'; ForceText = $true },
    @{ Cell = 'E28'; Value = '  -0.35%  '; ForceText = $false },
    @{ Cell = 'D29'; Value = '0.999'; ForceText = $true },
    @{ Cell = 'E29'; Value = '  +0.04%  '; ForceText = $false },
    @{ Cell = 'E30'; Value = '  -3.47%  '; ForceText = $false },
    @{ Cell = 'D31'; Value = '2.12'; ForceText = $true },
    @{ Cell = 'E31'; Value = '  -0.85%  '; ForceText = $false },
    @{ Cell = 'D32'; Value = '0.999'; ForceText = $true },
    @{ Cell = 'E32'; Value = '  -0.05%  '; ForceText = $false },
    @{ Cell = 'E33'; Value = '  -3.85%  '; ForceText = $false },
    @{ Cell = 'D34'; Value = '23.02'; ForceText = $true },
    @{ Cell = 'E34'; Value = '  -2.86%  '; ForceText = $false },
    @{ Cell = 'D35'; Value = '5.46'; ForceText = $true },
    @{ Cell = 'E35'; Value = '  +2.34%  '; ForceText = $false },
    @{ Cell = 'D36'; Value = '1.60'; ForceText = $true },
    @{ Cell = 'E36'; Value = '  +2.29%  '; ForceText = $false },
    @{ Cell = 'E37'; Value = '  -0.36%  '; ForceText = $false },
    @{ Cell = 'D38'; Value = '6.88'; ForceText = $true },
    @{ Cell = 'E38'; Value = '  -3.22%  '; ForceText = $false },
    @{ Cell = 'D39'; Value = '30.95'; ForceText = $true },
    @{ Cell = 'E39'; Value = '  -3.52%  '; ForceText = $false },
    @{ Cell = 'D40'; Value = '3.445.22'; ForceText = $false },
    @{ Cell = 'E40'; Value = '  -1.69%  '; ForceText = $false },
    @{ Cell = 'D41'; Value = '0.0772'; ForceText = $true },
    @{ Cell = 'E41'; Value = '  +1.26%  '; ForceText = $false },
    @{ Cell = 'D42'; Value = '42.49'; ForceText = $true },
    @{ Cell = 'E42'; Value = '  +0.38%  '; ForceText = $false },
    @{ Cell = 'E43'; Value = '  -2.95%  '; ForceText = $false },
    @{ Cell = 'E44'; Value = '  -1.46%  '; ForceText = $false },
    @{ Cell = 'D45'; Value = '1.66'; ForceText = $true },
    @{ Cell = 'E45'; Value = '  -3.70%  '; ForceText = $false },
    @{ Cell = 'D46'; Value = '1.16'; ForceText = $true },
    @{ Cell = 'E46'; Value = '  -5.49%  '; ForceText = $false },
    @{ Cell = 'D47'; Value = '2.536.48'; ForceText = $false },
    @{ Cell = 'E47'; Value = '  -3.45%  '; ForceText = $false },
    @{ Cell = 'D48'; Value = '6.88'; ForceText = $true },
    @{ Cell = 'E48'; Value = '  +1.53%  '; ForceText = $false },
    @{ Cell = 'D49'; Value = '22.53'; ForceText = $true },
    @{ Cell = 'E49'; Value = '  -2.82%  '; ForceText = $false },
    @{ Cell = 'B50'; Value = 'FirstDigitalUSD'; ForceText = $false },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'; ForceText = $false },
    @{ Cell = 'D50'; Value = '1.00'; ForceText = $true },
    @{ Cell = 'E50'; Value = '  +0.08%  '; ForceText = $false },
    @{ Cell = 'B51'; Value = 'dogwifhat'; ForceText = $false },
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'; ForceText = $false },
    @{ Cell = 'D51'; Value = '2.18'; ForceText = $true },
    @{ Cell = 'E51'; Value = '  -4.80%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $value = $u.Value
    if ($u.ForceText) {
        $value = "'" + $value
    }
    $ws.Range($u.Cell).Value = $value
}
